# ECUSpecs.xlsx update:
#  - insert a new "Heating Load (KW)" column after "Cooling Load (KW)" (old col E -> new col F, etc.)
#  - update BTU / cooling-load figures for the 60K, 36K and 18K units
#  - populate the new Heating Load (KW) column
#  - add/extend cell comments to document the new column and the "@ 125 degrees" metric note
#  - move the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new column (shifts old E:L -> F:M, formulas + comments anchors
#    for the cell VALUES move automatically; comments themselves do not, so
#    those are repaired explicitly in step 4 below).
# ---------------------------------------------------------------------------
$ws.Columns("E").Insert()

# Give the newly inserted column the same width as column D (17.85546875 in
# Excel's stored character-width units -> ColumnWidth 17 reproduces that).
$ws.Columns("E").ColumnWidth = 17

# ---------------------------------------------------------------------------
# 2. Header for the new column.
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "Heating Load (KW)"

# ---------------------------------------------------------------------------
# 3. Updated source data (cols B:D) + new Heating Load (KW) data (col E).
# ---------------------------------------------------------------------------
# Row 2 - 60K
$ws.Range("B2").Value = 54200
$ws.Range("C2").Value = 37100
$ws.Range("D2").Value = 12.231999999999999
$ws.Range("E2").Value = 12.646000000000001

# Row 3 - 36K
$ws.Range("B3").Value = 34100
$ws.Range("C3").Value = 31000
$ws.Range("D3").Value = 7.0380000000000003
$ws.Range("E3").Value = 10.371

# Row 4 - 18K
$ws.Range("B4").Value = 18500
$ws.Range("C4").Value = 14300
$ws.Range("D4").Value = 4.9050000000000002
$ws.Range("E4").Value = 4.9560000000000004

# Row 5 - HDT (B:D unchanged, only the new Heating Load value is added)
$ws.Range("E5").Value = 11.7

# ---------------------------------------------------------------------------
# 4. Repair / extend comments.
#    Column insert does not relocate existing comments, so the comment that
#    used to sit over the old column E ("Cost") etc. is still anchored to
#    E1..L1 even though the underlying header text moved to F1..M1. Fix this
#    by rewriting the text of every existing comment in place (so authorship
#    / formatting is preserved) to match the text that belongs at that
#    address now, then add one brand-new comment for the new M1 header.
# ---------------------------------------------------------------------------
$ws.Range("B1").Comment.Text("Clary Capt Deryk L:`nThe maximum BTU/hr capacity`nUsing ""@ 125 degrees"" metrics.")
$ws.Range("D1").Comment.Text("Clary Capt Deryk L:`nThe maximum amount of electrical load produced by the unit.`nUsing ""@ 125 degrees"" metrics.")

$ws.Range("E1").Comment.Text("Clary Capt Deryk L:`nPower draw when heating in kilowatts.")
$ws.Range("F1").Comment.Text("Clary Capt Deryk L:`nThe cost of an individual unit. Current values obtained from TDM Catalyst.`n`nThe current inputs are solely the procurement cost of a single unit. User can change this to include operations and maintence cost if the data is available.")
$ws.Range("G1").Comment.Text("Clary Capt Deryk L:`nThe weight of the unit.")
$ws.Range("H1").Comment.Text("Clary Capt Deryk L:`nThe amount of space taken up by the unit.")
$ws.Range("I1").Comment.Text("Clary Capt Deryk L:`nTrue: If the unit is window-mounted.`nFalse: If the unit is not window-mounted.")
$ws.Range("J1").Comment.Text("Clary Capt Deryk L:`nFormula calculated: The cost per unit divided by the cooling capacity.`n`nThese formula-calculated columns are not needed but are good metrics to fact-check the output of the tool.")
$ws.Range("K1").Comment.Text("Clary Capt Deryk L:`nFormula calculated: The amount of power used per BTU/hr")
$ws.Range("L1").Comment.Text("Clary Capt Deryk L:`nFormula calculated: The amount of space taken up per BTU/hr")

$ws.Range("M1").AddComment("Clary Capt Deryk L:`nFormula calculated: The number of pounds per BTU/hr")

# ---------------------------------------------------------------------------
# 5. Selection moves to F5.
# ---------------------------------------------------------------------------
$ws.Range("F5").Select()
